$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("B2").Value = 2000
$ws.Range("B3").Value = 467

# Update the active selection / view
$ws.Range("B4").Select()
